# "Use 16 time slices"
#
# The "Region-Time Slices" sheet defines the time-slice tree via three
# small lookup tables (season_info, weekly_info, daynite_info) plus a
# per-region summary (columns E:G). Previously every month/weekday/hour
# mapped to a single generic slice "A" (i.e. 1 time slice total). This
# change splits them into 4 seasons (Wi/Sp/Su/Au) x 1 weekly type (Day) x
# 4 day/night periods (N/D/P/E) = 16 time slices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Region-Time Slices")

# season_info table (K4:L16) - Month -> Season
$seasons = @("Wi","Wi","Sp","Sp","Sp","Su","Su","Su","Au","Au","Au","Wi")
for ($i = 0; $i -lt $seasons.Length; $i++) {
    $ws.Range("L" + (5 + $i)).Value = $seasons[$i]
}

# weekly_info table (N4:O11) - Day -> Type
for ($i = 0; $i -lt 7; $i++) {
    $ws.Range("O" + (5 + $i)).Value = "Day"
}

# daynite_info table (Q4:R28) - Hour -> Type
$daynite = @("N","N","N","N","N","N","N","D","D","D","D","D","D","D","D","D","D","P","P","P","E","E","E","E")
for ($i = 0; $i -lt $daynite.Length; $i++) {
    $ws.Range("R" + (5 + $i)).Value = $daynite[$i]
}

# Per-region summary columns (E:G) for the first region block (IE / National)
$ws.Range("E5").Value = "Wi"
$ws.Range("F5").Value = "Day"
$ws.Range("G5").Value = "N"

$ws.Range("E6").Value = "Sp"
$ws.Range("G6").Value = "D"

$ws.Range("E7").Value = "Su"
$ws.Range("G7").Value = "P"

$ws.Range("E8").Value = "Au"
$ws.Range("G8").Value = "E"
